$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filename string used across F2:F4 (shared string "p2_true_table.json" -> "documents_true_table.json")
$ws.Range("F2:F4").Value = "documents_true_table.json"

# Row 2 (§ 275.0-2_P2)
$ws.Range("B2").Value = 31
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 22

# Row 3 (§ 275.0-5_P2)
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 3

# Row 4 (§ 275.0-7_P2)
$ws.Range("B4").Value = 31
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 23
